# "correction to the best Val R^2"
#
# - Rebuild the "opt" sheet with the corrected Final Summary Metrics table
#   (adds a "Best Val R^2" row at 0.32011, the corrected figure) plus the
#   new Optimized Architecture Strategy + Training & Performance Benchmark
#   sections.
# - Re-point the active sheet/selection back to "exp" (with its own
#   selection moved to H18), and drop the "opt" sheet's tabSelected flag,
#   moving its saved selection to B15.

$wb = $excel.ActiveWorkbook

$expSheet = $wb.Worksheets.Item(1)
$optSheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Rebuild "opt" sheet content from scratch.
# ---------------------------------------------------------------------
$optSheet.Cells.UnMerge()
$optSheet.Cells.Clear()

# All data cells in the rebuilt tables are centered (style "2" in the
# original workbook); apply that per populated block (leaving the blank
# separator rows 8/16 and the unused B/C columns of the benchmark block
# untouched) and then layer bold-centered ("3"/"4") onto the header rows.
$optSheet.Range("A1:C7").HorizontalAlignment = -4108
$optSheet.Range("A9:C15").HorizontalAlignment = -4108
$optSheet.Range("A17:A24").HorizontalAlignment = -4108

# --- Section 1: Final Summary Metrics -------------------------------
$optSheet.Range("A1:C1").Font.Bold = $true
$optSheet.Range("A1").Value2 = "Final Summary Metrics"
$optSheet.Range("A1:C1").Merge()

$optSheet.Range("A2:C2").Font.Bold = $true
$optSheet.Range("A2").Value2 = "Metric"
$optSheet.Range("B2").Value2 = "Value"
$optSheet.Range("C2").Value2 = "Significance"

$optSheet.Range("A3").Value2 = "Avg Train Loss"
$optSheet.Range("B3").Value2 = 0.15376
$optSheet.Range("C3").Value2 = "Robust training convergence."

$optSheet.Range("A4").Value2 = "Avg Val Loss"
$optSheet.Range("B4").Value2 = 0.12544
$optSheet.Range("C4").Value2 = "Efficient error minimization on validation sets."

$optSheet.Range("A5").Value2 = "Avg Train R^2"
$optSheet.Range("B5").Value2 = 0.15202
$optSheet.Range("C5").Value2 = "Stable baseline training performance."

$optSheet.Range("A6").Value2 = "Avg Val R^2"
$optSheet.Range("B6").Value2 = 0.03891
$optSheet.Range("C6").Value2 = "Positive average generalization across all folds."

$optSheet.Range("A7").Value2 = "Best Val R^2"
$optSheet.Range("B7").Value2 = 0.32011
$optSheet.Range("C7").Value2 = "Peak Potential: The highest variance explained in a single run."
$optSheet.Range("C7").Characters(16, 50).Font.Bold = $false

# --- Section 2: Optimized Architecture Strategy ----------------------
$optSheet.Range("A9:C9").Font.Bold = $true
$optSheet.Range("A9").Value2 = "Optimized Architecture Strategy"
$optSheet.Range("A9:C9").Merge()

$optSheet.Range("A10:C10").Font.Bold = $true
$optSheet.Range("A10").Value2 = "Component"
$optSheet.Range("B10").Value2 = "Optimal Selection"
$optSheet.Range("C10").Value2 = "Reasoning from Experiments"

$optSheet.Range("A11").Value2 = "Kernel Size"
$optSheet.Range("B11").Value2 = "3x3"
$optSheet.Range("C11").Value2 = "Top Performer: Achieved the highest individual R^2 potential (0.238–0.320)."

$optSheet.Range("A12").Value2 = "Depth"
$optSheet.Range("B12").Value2 = "3 Conv Blocks"
$optSheet.Range("C12").Value2 = "Best Balance: Maintained the lowest validation loss (0.120) while capturing complex features."
$optSheet.Range("C12").Characters(14, 77).Font.Bold = $false

$optSheet.Range("A13").Value2 = "Regularization"
$optSheet.Range("B13").Value2 = "Dropout 0.6"
$optSheet.Range("C13").Value2 = "Critical: High dropout was essential to stop the model from overfitting on small biomass samples."
$optSheet.Range("C13").Characters(10, 90).Font.Bold = $false

$optSheet.Range("A14").Value2 = "Optimization"
$optSheet.Range("B14").Value2 = "Adam + 0.0 WD"
$optSheet.Range("C14").Value2 = "Winner: Weight Decay of 0.0 combined with high dropout yielded the best validation R^2."

$optSheet.Range("A15").Value2 = "Augmentation"
$optSheet.Range("B15").Value2 = "Light/Medium"
$optSheet.Range("C15").Value2 = "Stability: Required to maintain positive R^2 values and reduce validation variance."

# --- Section 3: Training & Performance Benchmark ---------------------
$optSheet.Range("A17").Font.Bold = $true
$optSheet.Range("A17").Value2 = "Training & Performance Benchmark"

$optSheet.Range("A18").Value2 = "Metric,Optimized Baseline Value"
$optSheet.Range("A19").Value2 = "Number of Parameters,~1.5M (approx. for 256 filters)"
$optSheet.Range("A20").Value2 = "Avg Training Loss,0.15376"
$optSheet.Range("A21").Value2 = "Avg Validation Loss,0.12544"
$optSheet.Range("A22").Value2 = "Final Train R2,0.15202"
$optSheet.Range("A23").Value2 = "Final Val R2,0.03891"
$optSheet.Range("A24").Value2 = "Max R2 Achieved,0.32011"

# ---------------------------------------------------------------------
# 2. Fix up sheet selections / active sheet. "opt" is set first so that
#    activating "exp" afterwards leaves it as the final active sheet
#    (selecting a range on a sheet re-activates that sheet).
# ---------------------------------------------------------------------
$optSheet.Activate()
$optSheet.Range("B15").Select()

$expSheet.Activate()
$expSheet.Range("H18").Select()
